{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst targetText = \"M\u00c1QUINASEscola PRO-TEC\";\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetText) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Anchor paragraph not found: \" + targetText);\n}\n\n// The three paragraphs immediately following the anchor are removed:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 ...\" copyright paragraph\nconst toDelete = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3 && i < items.length; i++) {\n  toDelete.push(items[i]);\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$target = \"M\u00c1QUINASEscola PRO-TEC\"\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq $target) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph not found: $target\"\n}\n\n# Delete the three paragraphs immediately following the anchor paragraph:\n#   1) an empty paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) the \"\u00a9 2020 ...\" copyright paragraph\nfor ($n = 1; $n -le 3; $n++) {\n    $p = $d.Paragraphs.Item($anchorIndex + 1)\n    $p.Range.Delete()\n}\n"}
